$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.862.36"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").Value = "2.089.63"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.75"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.32"

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("E10").Value = "  +2.58%  "

$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("D12").Value = "2.396.81"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.73"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.43"
$ws.Range("E14").Value = "  +3.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.772"
$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("E16").Value = "  +2.60%  "

$ws.Range("D17").Value = "2.088.43"
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").Value = "37.784.68"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.90"
$ws.Range("E20").Value = "  +3.16%  "

$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.05"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("E24").Value = "  -0.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.85"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("E27").Value = "  +9.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.09"
$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.60"
$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("E31").Value = "  +2.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.71"
$ws.Range("E32").Value = "  +3.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").Value = "  +2.43%  "

$ws.Range("E34").Value = "  +3.43%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.51"
$ws.Range("E36").Value = "  +6.13%  "

$ws.Range("E37").Value = "  +2.69%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("E40").Value = "  +2.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.74"
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("E43").Value = "  +1.53%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.463.40"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.35"
$ws.Range("E45").Value = "  +8.18%  "

$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.25"
$ws.Range("E47").Value = "  +5.96%  "

$ws.Range("E48").Value = "  +5.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.49"
$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.78"
$ws.Range("E51").Value = "  +6.35%  "
